$d = $word.ActiveDocument

$replacements = @(
    @("2025-08-14 Thursday", "2025-08-15 Friday"),
    @("617÷6=", "651÷6="),
    @("642÷9=", "453÷7="),
    @("669÷8=", "434÷3="),
    @("439÷5=", "486÷7="),
    @("777÷3=", "989÷6="),
    @("485÷7=", "441÷4="),
    @("389÷3=", "578÷7="),
    @("540÷8=", "818÷3="),
    @("265÷6=", "682÷2="),
    @("778÷5=", "891÷5="),
    @("722÷9=", "358÷5="),
    @("689÷9=", "595÷4="),
    @("621÷9=", "839÷8="),
    @("793÷9=", "124÷5="),
    @("281÷4=", "710÷5="),
    @("838÷3=", "956÷5="),
    @("971÷7=", "858÷3="),
    @("247÷5=", "114÷7="),
    @("419÷5=", "470÷2="),
    @("842÷9=", "574÷4="),
    @("766÷3=", "721÷3="),
    @("497÷2=", "236÷3="),
    @("465÷8=", "292÷6="),
    @("702÷7=", "354÷4="),
    @("902÷2=", "524÷5=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
